# This script applies a cyclic re-shuffle of the data held in rows 11, 12
# and 13 of the "Artfynd" sheet:
#   - the record that used to live in row 13 now lives in row 11
#   - the record that used to live in row 11 now lives in row 12
#   - the record that used to live in row 12 now lives in row 13
#
# Rather than moving whole rows (which would risk Excel re-interpreting
# text such as the date strings in columns Y/AA as real dates), only the
# individual cells whose contents actually change are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (becomes the former row 13 record) ---------------------------
$ws.Range("A11").Value2 = 131046735
$ws.Range("B11").Value2 = 57884
$ws.Range("D11").Value2 = "NT"
$ws.Range("E11").Value2 = 100109
$ws.Range("F11").Value2 = "Tretåig hackspett"
$ws.Range("G11").Value2 = "Picoides tridactylus"
$ws.Range("H11").Value2 = "(Linnaeus, 1758)"
$ws.Range("M11").Value2 = "nyligen använt bo"
$ws.Range("Q11").Value2 = 402448
$ws.Range("R11").Value2 = 6818295
$ws.Range("Z11").Value2 = "16:54"
$ws.Range("AB11").Value2 = "16:54"

# --- Row 12 (becomes the former row 11 record) ---------------------------
$ws.Range("A12").Value2 = 131046763
$ws.Range("B12").Value2 = 92267
$ws.Range("D12").Value2 = "VU"
$ws.Range("E12").Value2 = 1209
$ws.Range("F12").Value2 = "Rynkskinn"
$ws.Range("G12").Value2 = "Hermanssonia centrifuga"
$ws.Range("H12").Value2 = "(P. Karst.) Zmitr."
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value2 = 402378
$ws.Range("R12").Value2 = 6818392
$ws.Range("Z12").Value2 = "17:02"
$ws.Range("AB12").Value2 = "17:02"
$ws.Range("AC12").ClearContents()

# --- Row 13 (becomes the former row 12 record) ---------------------------
$ws.Range("A13").Value2 = 131046788
$ws.Range("B13").Value2 = 57884
$ws.Range("D13").Value2 = "NT"
$ws.Range("E13").Value2 = 100109
$ws.Range("F13").Value2 = "Tretåig hackspett"
$ws.Range("G13").Value2 = "Picoides tridactylus"
$ws.Range("H13").Value2 = "(Linnaeus, 1758)"
$ws.Range("M13").Value2 = "färska spår"
$ws.Range("Q13").Value2 = 402473
$ws.Range("R13").Value2 = 6818425
$ws.Range("Z13").Value2 = "16:47"
$ws.Range("AB13").Value2 = "16:47"
$ws.Range("AC13").Value2 = "Färska ringhack (gran)"
